# Generate Report for Handoff
#
# Localization status moved from "In Translation" to "Ready for handoff",
# with a fresh handoff timestamp recorded on the Overview sheet and on
# each per-locale sheet (zh-cn, de-de). Also widen the "Status"-related
# columns that now hold the longer "Ready for handoff" label, matching
# the auto-fit width Excel would have produced for the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# --- Overview sheet: per-locale status + latest handoff generate date ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-18 06:37:41"

# --- zh-cn sheet: status + latest handoff datetime ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-18 06:37:36"

# --- de-de sheet: status + latest handoff datetime ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-18 06:37:41"

# --- Widen the Status columns to fit "Ready for handoff" ---
# target stored OOXML column width is 17.2159881591797; the host quantizes
# ColumnWidth writes to 1/6-character steps, so 49/3 is the closest input
# that round-trips to the nearest representable width (17.1666...).
$overview.Columns.Item(5).ColumnWidth = 49/3
$overview.Columns.Item(6).ColumnWidth = 49/3
$zhcn.Columns.Item(3).ColumnWidth = 49/3
$dede.Columns.Item(3).ColumnWidth = 49/3
